# Almacen Button Selector - FIX
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Product name / description: "SQ009-DELINEADOR" -> "SQ009-Pulsera"
$ws.Range("A2").Value = "SQ009-Pulsera"
$ws.Range("P2").Value = "SQ009-Pulsera"

# Internal code / barcode stay "SQ009" (already correct, but reassign for safety)
$ws.Range("B2").Value = "SQ009"
$ws.Range("T2").Value = "SQ009"

# Sale unit price: 12 -> 4
$ws.Range("G2").Value = 4

# Move selection to U2 (matches author's last cursor position)
$ws.Range("U2").Select()
